# Apply crypto price/volume updates per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.723.74'
$ws.Range('E2').Value = '  -3.80%  '
$ws.Range('D3').Value = '3.051.05'
$ws.Range('E3').Value = '  -3.09%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'544.64"
$ws.Range('E5').Value = '  -4.27%  '
$ws.Range('D6').Value = "'134.02"
$ws.Range('E6').Value = '  -10.59%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.045.82'
$ws.Range('E8').Value = '  -3.02%  '
$ws.Range('D9').Value = "'0.489"
$ws.Range('E9').Value = '  -2.93%  '
$ws.Range('E10').Value = '  -3.96%  '
$ws.Range('D11').Value = "'6.38"
$ws.Range('E11').Value = '  -11.05%  '
$ws.Range('D12').Value = "'0.459"
$ws.Range('E12').Value = '  -2.43%  '
$ws.Range('D13').Value = "'34.66"
$ws.Range('E13').Value = '  -4.93%  '
$ws.Range('D14').Value = "'0.0000214"
$ws.Range('E14').Value = '  -5.45%  '
$ws.Range('D15').Value = '3.540.66'
$ws.Range('E15').Value = '  -3.22%  '
$ws.Range('D16').Value = '62.729.92'
$ws.Range('E16').Value = '  -3.83%  '
$ws.Range('E17').Value = '  -2.58%  '
$ws.Range('D18').Value = '3.042.87'
$ws.Range('E18').Value = '  -3.31%  '
$ws.Range('D19').Value = "'6.62"
$ws.Range('E19').Value = '  -3.28%  '
$ws.Range('D20').Value = "'480.37"
$ws.Range('E20').Value = '  -11.10%  '
$ws.Range('D21').Value = "'13.34"
$ws.Range('E21').Value = '  -5.10%  '
$ws.Range('D22').Value = "'0.697"
$ws.Range('E22').Value = '  -2.30%  '
$ws.Range('D23').Value = "'7.05"
$ws.Range('E23').Value = '  -6.17%  '
$ws.Range('D24').Value = "'77.17"
$ws.Range('E24').Value = '  -2.85%  '
$ws.Range('D25').Value = "'12.15"
$ws.Range('E25').Value = '  -7.09%  '
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').Value = "'8.26"
$ws.Range('E27').Value = '  -8.96%  '
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').Value = "'2.69"
$ws.Range('E28').Value = '  -4.76%  '
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('D30').Value = "'1.93"
$ws.Range('E30').Value = '  -11.09%  '
$ws.Range('D31').Value = "'26.13"
$ws.Range('E31').Value = '  -1.42%  '
$ws.Range('E32').Value = '  -2.86%  '
$ws.Range('D33').Value = "'2.48"
$ws.Range('E33').Value = '  -7.67%  '
$ws.Range('D34').Value = "'59.06"
$ws.Range('E34').Value = '  +11.57%  '
$ws.Range('D35').Value = "'508.83"
$ws.Range('E35').Value = '  -8.19%  '
$ws.Range('D36').Value = "'5.94"
$ws.Range('E36').Value = '  -4.15%  '
$ws.Range('D37').Value = "'5.06"
$ws.Range('E37').Value = '  -8.36%  '
$ws.Range('D38').Value = "'0.0397"
$ws.Range('E38').Value = '  -11.99%  '
$ws.Range('D39').Value = '3.083.18'
$ws.Range('E39').Value = '  +0.34%  '
$ws.Range('D40').Value = "'0.0786"
$ws.Range('E40').Value = '  -5.52%  '
$ws.Range('D41').Value = "'0.118"
$ws.Range('E41').Value = '  -4.79%  '
$ws.Range('D42').Value = "'8.01"
$ws.Range('E42').Value = '  -4.34%  '
$ws.Range('D43').Value = "'2.59"
$ws.Range('E43').Value = '  -11.80%  '
$ws.Range('D44').Value = "'0.251"
$ws.Range('E44').Value = '  -4.43%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').Value = "'2.02"
$ws.Range('E46').Value = '  -9.57%  '
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').Value = "'119.52"
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = "'24.34"
$ws.Range('E48').Value = '  -4.15%  '
$ws.Range('D49').Value = "'0.107"
$ws.Range('E49').Value = '  -3.51%  '
$ws.Range('B50').Value = 'CoreDAO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D50').Value = "'2.38"
$ws.Range('E50').Value = '  +59.77%  '
$ws.Range('B51').Value = 'PEPE'
$ws.Range('C51').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D51').Value = '0.0₃0492'
$ws.Range('E51').Value = '  -7.81%  '
